$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 4
$ws.Range("D12").Value = "finished and pushed student major/minor, need to ask about other options… Now it’s time to find and fix bugs of empty columns…"
$ws.Range("D12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 42.75
$ws.Range("D13").Select()
